# Update "想去人数" (F column) figures for two events that appear on both
# the "展览" sheet and the aggregated "全部类型" sheet.
#   - "安徽·崩坏同人only 爱莉希雅同人生日会"          20   -> 21
#   - "合肥·第九届环形宇宙动漫游戏嘉年华"              3708 -> 3716
#   - "合肥·心动恋章·冬日序国乙&代号鹄同人only"        193  -> 196

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 21
$ws1.Range("F4").Value = 3716
$ws1.Range("F7").Value = 196

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 21
$ws4.Range("F8").Value = 3716
$ws4.Range("F12").Value = 196
